$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation,
# then restore the cell's original style so no stray number format is left behind.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "30.284.78"
$ws.Range("E2").Value = "  -0.13%  "

$ws.Range("D3").Value = "1.914.12"
$ws.Range("E3").Value = "  -1.09%  "

Set-TextValue "D4" "0.9999"
$ws.Range("E4").Value = "  +0.06%  "

Set-TextValue "D5" "0.7404"
$ws.Range("E5").Value = "  -3.00%  "

Set-TextValue "D6" "243.99"
$ws.Range("E6").Value = "  -1.81%  "

Set-TextValue "D7" "0.9999"
$ws.Range("E7").Value = "  +0.15%  "

Set-TextValue "D8" "0.3149"
$ws.Range("E8").Value = "  -2.19%  "

Set-TextValue "D9" "27.24"
$ws.Range("E9").Value = "  -3.78%  "

Set-TextValue "D10" "0.07011"
$ws.Range("E10").Value = "  -1.53%  "

Set-TextValue "D11" "0.7834"
$ws.Range("E11").Value = "  -0.97%  "

Set-TextValue "D12" "0.07967"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").Value = "1.904.78"
$ws.Range("E13").Value = "  -1.60%  "

Set-TextValue "D14" "5.310"
$ws.Range("E14").Value = "  -1.34%  "

Set-TextValue "D15" "92.12"
$ws.Range("E15").Value = "  -2.81%  "

Set-TextValue "D16" "14.39"
$ws.Range("E16").Value = "  -2.41%  "

$ws.Range("D17").Value = "30.253.48"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D18" "245.64"
$ws.Range("E18").Value = "  -3.46%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D19" "5.849"
$ws.Range("E19").Value = "  +0.78%  "

Set-TextValue "D20" "0.000007868"
$ws.Range("E20").Value = "  -2.11%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.168.10"
$ws.Range("E21").Value = "  -0.87%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D22" "0.9999"
$ws.Range("E22").Value = "  +0.24%  "

$ws.Range("E23").Value = "  +0.02%  "

Set-TextValue "D24" "6.677"
$ws.Range("E24").Value = "  -2.32%  "

Set-TextValue "D25" "9.494"
$ws.Range("E25").Value = "  -1.05%  "

Set-TextValue "D26" "164.93"
$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("E27").Value = "  -0.52%  "

Set-TextValue "D28" "0.1275"
$ws.Range("E28").Value = "  -6.02%  "

Set-TextValue "D29" "2.126"
$ws.Range("E29").Value = "  -8.40%  "

Set-TextValue "D30" "1.348"
$ws.Range("E30").Value = "  -1.85%  "

Set-TextValue "D31" "1.548"
$ws.Range("E31").Value = "  +1.31%  "

Set-TextValue "D32" "4.340"
$ws.Range("E32").Value = "  -2.31%  "

Set-TextValue "D33" "4.092"
$ws.Range("E33").Value = "  -1.44%  "

Set-TextValue "D34" "0.05214"
$ws.Range("E34").Value = "  +0.39%  "

Set-TextValue "D35" "1.312"
$ws.Range("E35").Value = "  +1.46%  "

Set-TextValue "D36" "0.7523"
$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("E37").Value = "  -0.47%  "

Set-TextValue "D38" "0.01951"
$ws.Range("E38").Value = "  -1.16%  "

Set-TextValue "D39" "2.801"
$ws.Range("E39").Value = "  -0.02%  "

Set-TextValue "D40" "6.392"
$ws.Range("E40").Value = "  -0.94%  "

Set-TextValue "D41" "76.07"
$ws.Range("E41").Value = "  -3.17%  "

Set-TextValue "D42" "0.4509"
$ws.Range("E42").Value = "  -0.46%  "

Set-TextValue "D43" "1.949"
$ws.Range("E43").Value = "  -2.49%  "

Set-TextValue "D44" "0.9993"
$ws.Range("E44").Value = "  +0.05%  "

$ws.Range("E45").Value = "  +2.80%  "

Set-TextValue "D46" "0.8335"
$ws.Range("E46").Value = "  -0.51%  "

Set-TextValue "D47" "9.900"
$ws.Range("E47").Value = "  +0.58%  "

Set-TextValue "D48" "101.22"
$ws.Range("E48").Value = "  -1.13%  "

Set-TextValue "D49" "37.43"
$ws.Range("E49").Value = "  +0.12%  "

Set-TextValue "D50" "0.1215"
$ws.Range("E50").Value = "  +1.65%  "

Set-TextValue "D51" "941.92"
$ws.Range("E51").Value = "  -4.76%  "

